$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $helper = $ws.Range("Z100")
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy()
    $cell.PasteSpecial(-4163)
    $helper.Clear()
    $excel.CutCopyMode = 0
}

$ws.Range("D2").Value = '29.878.03'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '1.872.55'
$ws.Range("E3").Value = '  -1.17%  '
Set-TextValue $ws.Range("D4") '0.9998'
$ws.Range("E4").Value = '  -0.05%  '
Set-TextValue $ws.Range("D5") '0.7350'
$ws.Range("E5").Value = '  -5.02%  '
Set-TextValue $ws.Range("D6") '242.03'
$ws.Range("E6").Value = '  -1.10%  '
Set-TextValue $ws.Range("D7") '0.9998'
$ws.Range("E7").Value = '  -0.04%  '
Set-TextValue $ws.Range("D8") '0.3154'
$ws.Range("E8").Value = '  +0.48%  '
Set-TextValue $ws.Range("D9") '24.66'
$ws.Range("E9").Value = '  -4.48%  '
Set-TextValue $ws.Range("D10") '0.07095'
$ws.Range("E10").Value = '  -2.38%  '
Set-TextValue $ws.Range("D11") '0.08461'
$ws.Range("E11").Value = '  +0.67%  '
Set-TextValue $ws.Range("D12") '0.7510'
$ws.Range("E12").Value = '  -2.96%  '
$ws.Range("D13").Value = '1.881.02'
$ws.Range("E13").Value = '  +0.60%  '
Set-TextValue $ws.Range("D14") '5.372'
$ws.Range("E14").Value = '  -1.95%  '
Set-TextValue $ws.Range("D15") '92.55'
$ws.Range("E15").Value = '  -2.63%  '
$ws.Range("D16").Value = '29.881.64'
$ws.Range("E16").Value = '  +0.04%  '
Set-TextValue $ws.Range("D17") '6.040'
$ws.Range("E17").Value = '  -2.51%  '
Set-TextValue $ws.Range("D18") '13.59'
$ws.Range("E18").Value = '  -2.97%  '
Set-TextValue $ws.Range("D19") '243.23'
$ws.Range("E19").Value = '  -1.66%  '
Set-TextValue $ws.Range("D20") '0.000007809'
$ws.Range("E20").Value = '  -0.84%  '
Set-TextValue $ws.Range("D21") '0.9993'
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").Value = '2.122.14'
$ws.Range("E22").Value = '  +0.72%  '
Set-TextValue $ws.Range("D23") '7.920'
$ws.Range("E23").Value = '  -2.55%  '
Set-TextValue $ws.Range("D24") '0.9989'
$ws.Range("E24").Value = '  -0.15%  '
Set-TextValue $ws.Range("D25") '0.1563'
$ws.Range("E25").Value = '  -2.15%  '
Set-TextValue $ws.Range("D26") '9.326'
$ws.Range("E26").Value = '  -2.35%  '
Set-TextValue $ws.Range("D27") '163.99'
$ws.Range("E27").Value = '  +0.89%  '
Set-TextValue $ws.Range("D28") '18.61'
$ws.Range("E28").Value = '  -0.99%  '
Set-TextValue $ws.Range("D29") '2.023'
$ws.Range("E29").Value = '  -0.86%  '
Set-TextValue $ws.Range("D30") '1.466'
$ws.Range("E30").Value = '  +3.28%  '
Set-TextValue $ws.Range("D31") '4.552'
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("E32").Value = '  -1.38%  '
Set-TextValue $ws.Range("D33") '4.278'
$ws.Range("E33").Value = '  +4.10%  '
Set-TextValue $ws.Range("D34") '0.05332'
$ws.Range("E34").Value = '  -2.61%  '
Set-TextValue $ws.Range("D35") '1.235'
$ws.Range("E35").Value = '  -1.06%  '
Set-TextValue $ws.Range("D36") '0.7507'
$ws.Range("E36").Value = '  -0.11%  '
Set-TextValue $ws.Range("D37") '0.9995'
$ws.Range("E37").Value = '  -0.39%  '
Set-TextValue $ws.Range("D38") '2.701'
$ws.Range("E38").Value = '  +1.03%  '
Set-TextValue $ws.Range("D39") '0.01950'
$ws.Range("E39").Value = '  +0.67%  '
$ws.Range("E40").Value = '  -1.36%  '
Set-TextValue $ws.Range("D41") '0.4466'
$ws.Range("E41").Value = '  -0.65%  '
$ws.Range("D42").Value = '1.104.15'
$ws.Range("E42").Value = '  +0.76%  '
Set-TextValue $ws.Range("D43") '6.067'
$ws.Range("E43").Value = '  -0.24%  '
Set-TextValue $ws.Range("D44") '72.35'
$ws.Range("E44").Value = '  -2.64%  '
Set-TextValue $ws.Range("D45") '0.8652'
$ws.Range("E45").Value = '  +1.37%  '
$ws.Range("E46").Value = '  +0.09%  '
Set-TextValue $ws.Range("D47") '102.82'
$ws.Range("E47").Value = '  +0.16%  '
Set-TextValue $ws.Range("D48") '7.724'
$ws.Range("E48").Value = '  +1.53%  '
Set-TextValue $ws.Range("D49") '3.077'
$ws.Range("E49").Value = '  +2.15%  '
Set-TextValue $ws.Range("D50") '1.838'
$ws.Range("E50").Value = '  -2.89%  '
$ws.Range("D51").Value = '2.020.08'
$ws.Range("E51").Value = '  -0.66%  '
